$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows corresponding to additional daily dates (31-08-2021 .. 06-09-2021)
$newRows = @(
    @{ Row = 171; Date = "31-08-2021"; Values = @(37, 7, 9, -2, 0, -2, -5, -31, -19, 6) },
    @{ Row = 172; Date = "01-09-2021"; Values = @(38, 7, 9, -2, 0, -2, -12, -27, -17, 7) },
    @{ Row = 173; Date = "02-09-2021"; Values = @(38, 7, 9, -2, 0, -2, -9, -30, -18, 7) },
    @{ Row = 174; Date = "03-09-2021"; Values = @(38, 7, 9, -2, 0, -2, -8, -34, -14, 7) },
    @{ Row = 175; Date = "06-09-2021"; Values = @(38, 7, 9, -2, 0, -2, -8, -36, -12, 7) }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    # Column A holds the date as plain text (matches existing rows: no
    # explicit number format / style, stored as a shared string).
    $ws.Cells.Item($rowNum, 1).Value = $r.Date

    $colIndex = 2
    foreach ($v in $r.Values) {
        $ws.Cells.Item($rowNum, $colIndex).Value = $v
        $colIndex++
    }
}
